$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "데이터 사이언티스트의 미래 (1)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/data-scientists-future-1/#utm_source=rss&utm_medium=rss&utm_campaign=data-scientists-future-1"

$ws.Range("D29").Value = "[GAN 시리즈][DCLGAN] Dual Contrastive Learning for Unsupervised Image-to-Image Translation - 1편"
$ws.Range("E29").Value = "https://blog.promedius.ai/dclgan_1/"

$ws.Range("D32").Value = "텐서플로 함수와 그래프(tf.function, autograph)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/328"

$ws.Range("D52").Value = "사건까지 걸린 시간은?"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2330&utm_source=rss&utm_medium=rss&utm_campaign=%25ec%2582%25ac%25ea%25b1%25b4%25ea%25b9%258c%25ec%25a7%2580-%25ea%25b1%25b8%25eb%25a6%25b0-%25ec%258b%259c%25ea%25b0%2584%25ec%259d%2580"
